# Swap the NetServices and TrustSec_Devices SGACLs on the Matrix sheet.
#
# The Matrix sheet lists SGT rows (A column) against destination-SGT
# columns (D:K). Row 6 is "NetServices" and row 7 is "TrustSec_Devices".
# This swaps the SGACL values (Deny IP / Permit IP) assigned across the
# destination columns between those two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matrix")

# Grab the current values for the swap range (D:K) on both rows first,
# so the swap is not order-dependent.
$netServicesRange = $ws.Range("D6:K6")
$trustSecRange = $ws.Range("D7:K7")

$netServicesValues = $netServicesRange.Value2
$trustSecValues = $trustSecRange.Value2

$trustSecRange.Value2 = $netServicesValues
$netServicesRange.Value2 = $trustSecValues

$wb.Save()
